# Generate Report for Handoff
# Updates Priority and Latest Handoff/Generate Datetime values for the
# rows that were (re-)handed off since the last report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for rows 4-7
foreach ($r in 4..7) {
    $wsOverview.Range("G$r").Value = "2016-09-07 10:51:42"
}

# zh-cn sheet: Priority (column E) and Latest Handoff Datetime (column H) for rows 4-7
foreach ($r in 4..7) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-07 10:51:29"
}

# de-de sheet: Priority (column E) and Latest Handoff Datetime (column H) for rows 4-7
foreach ($r in 4..7) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-07 10:51:42"
}
